$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 3) of data to the kaspa_buys sheet, as recorded
# on 2025-08-31.
# Force column A to be treated as plain text so the date string is not
# auto-converted into a date serial number, then restore the default
# (unstyled) "Normal" style so the new cell matches the rest of the
# data rows, which carry no explicit cell style.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "08/31/2025"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 116.8280000000004
$ws.Range("C3").Value = 0.08559591878659194
$ws.Range("D3").Value = 10
